$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "21.05.202"
$ws.Range("B4").Value = "35,90€"
$ws.Range("C4").Value = "Karlsruhe "
$ws.Range("D4").Value = "Hannover "
$ws.Range("E4").Value = "Hannover "
